$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New "Exception Handling Demo" day-block (rows 135-140) ---------------
# Set the values/formula first (while the cells are still "fresh" new rows)
# so the calculation engine correctly evaluates the new SUM formula, then
# copy the formatting (borders / number formats / fonts / alignment) down
# from the previous day-block afterwards so the new rows reuse the existing
# cell styles instead of creating new ones.

# Row 135 is left as an empty separator row (no values).

# Row 136: new day entry
$ws.Range("A136").Value = 45692
$ws.Range("B136").Value = "Domm"
$ws.Range("D136").Value = 0.25

# Row 137: Meeting / Reconsile
$ws.Range("B137").Value = "Meeting"
$ws.Range("C137").Value = "Reconsile"
$ws.Range("D137").Value = 0

# Row 138: General Discussion
$ws.Range("C138").Value = "General Discussion"
$ws.Range("D138").Value = 0.25

# Row 139: Study / ASP.NET Core Final Demo
$ws.Range("B139").Value = "Study"
$ws.Range("C139").Value = "ASP.NET Core Final Demo"
$ws.Range("D139").Value = 7.5

# Row 140: Total row with formula summing the new block
$ws.Range("B140").Value = "Total"
$ws.Range("D140").Formula = "=SUM(D135:D139)"

# Copy the formatting of rows 124-129 (same blank/date/date/row/row/total
# shape) onto the new rows 135-140 so they pick up the existing styles.
$ws.Range("A124:D129").Copy() | Out-Null
$ws.Range("A135:D140").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# Unlike the other block's total rows, this total row has no formatted cell
# in column C, so clear it entirely after the format paste.
$ws.Range("C140").Clear() | Out-Null

# Update selection / scroll position to match the new extent of data.
$ws.Application.ActiveWindow.ScrollRow = 111
$ws.Range("A136:D141").Select()
